$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44614
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 21000
$ws.Range("P3").Value = 20500
$ws.Range("Q3").Value = '$/bandeja 18 kilos'
$ws.Range("S3").Value = 1139
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 44656
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 270
$ws.Range("N4").Value = 19000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 19500
$ws.Range("Q4").Value = '$/bandeja 18 kilos'
$ws.Range("S4").Value = 1083
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44291
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 17000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17500
$ws.Range("Q5").Value = '$/bandeja 18 kilos'
$ws.Range("S5").Value = 972
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 44263
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21500
$ws.Range("Q6").Value = '$/caja 18 kilos'
$ws.Range("S6").Value = 1194
$ws.Range("T6").Value = 18

# Row 7
$ws.Range("D7").Value = 44819
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 17000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 17500
$ws.Range("Q7").Value = '$/bandeja 10 kilos'
$ws.Range("S7").Value = 1750
$ws.Range("T7").Value = 10

# Row 8
$ws.Range("D8").Value = 44487
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("Q8").Value = '$/bandeja 10 kilos'
$ws.Range("S8").Value = 1450
$ws.Range("T8").Value = 10

# Row 9
$ws.Range("D9").Value = 44673
$ws.Range("L9").Value = 'Especial'
$ws.Range("M9").Value = 400
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14500
$ws.Range("Q9").Value = '$/bandeja 10 kilos'
$ws.Range("S9").Value = 1450
$ws.Range("T9").Value = 10

# Row 10
$ws.Range("D10").Value = 44784
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 300
$ws.Range("N10").Value = 19000
$ws.Range("O10").Value = 20000
$ws.Range("P10").Value = 19500
$ws.Range("Q10").Value = '$/bandeja 18 kilos'
$ws.Range("S10").Value = 1083
$ws.Range("T10").Value = 18

# Row 11
$ws.Range("D11").Value = 44489
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 26000
$ws.Range("O11").Value = 27000
$ws.Range("P11").Value = 26500
$ws.Range("Q11").Value = '$/bandeja 18 kilos'
$ws.Range("S11").Value = 1472
$ws.Range("T11").Value = 18

# Row 12
$ws.Range("D12").Value = 44616
$ws.Range("L12").Value = 'Segunda'
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 16000
$ws.Range("O12").Value = 17000
$ws.Range("P12").Value = 16500
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("S12").Value = 917
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 44602
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 270
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 21000
$ws.Range("P13").Value = 20500
$ws.Range("Q13").Value = '$/bandeja 18 kilos'
$ws.Range("S13").Value = 1139
$ws.Range("T13").Value = 18

# Row 14
$ws.Range("D14").Value = 44491
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 14000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 14500
$ws.Range("Q14").Value = '$/bandeja 10 kilos'
$ws.Range("S14").Value = 1450
$ws.Range("T14").Value = 10

# Row 15
$ws.Range("D15").Value = 44629
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 17000
$ws.Range("O15").Value = 18000
$ws.Range("P15").Value = 17500
$ws.Range("Q15").Value = '$/bandeja 18 kilos'
$ws.Range("S15").Value = 972
$ws.Range("T15").Value = 18

# Row 16
$ws.Range("D16").Value = 44307
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 250
$ws.Range("N16").Value = 19000
$ws.Range("O16").Value = 20000
$ws.Range("P16").Value = 19500
$ws.Range("Q16").Value = '$/bandeja 18 kilos'
$ws.Range("S16").Value = 1083
$ws.Range("T16").Value = 18

# Row 17
$ws.Range("D17").Value = 44418
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 240
$ws.Range("N17").Value = 10000
$ws.Range("O17").Value = 11000
$ws.Range("P17").Value = 10500
$ws.Range("Q17").Value = '$/bandeja 10 kilos'
$ws.Range("S17").Value = 1050
$ws.Range("T17").Value = 10

# Row 18
$ws.Range("D18").Value = 44789
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 250
$ws.Range("N18").Value = 19000
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 19500
$ws.Range("Q18").Value = '$/bandeja 18 kilos'
$ws.Range("S18").Value = 1083
$ws.Range("T18").Value = 18

# Row 19
$ws.Range("D19").Value = 44323
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 270
$ws.Range("N19").Value = 21000
$ws.Range("O19").Value = 22000
$ws.Range("P19").Value = 21500
$ws.Range("Q19").Value = '$/bandeja 18 kilos'
$ws.Range("S19").Value = 1194
$ws.Range("T19").Value = 18
